$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.575.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.596.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.61%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.990.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.585.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.893"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.543.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "267.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.61%  "

$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.78%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.60%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.90%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("E35").Value = "  -2.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0817"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.88%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.09%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.121"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.043.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.75%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.842.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.191"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.19%  "
